$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add offspring size (column G) and source URL (column H) under each declared animal

$ws.Range("G2").Value = "2000 eggs"
$ws.Range("H2").Value = "http://www.fcps.edu/islandcreekes/ecology/largemouth_bass.htm"

$ws.Range("G3").Value = "1 million"
$ws.Range("H3").Value = "http://www.gma.org/fogm/Gadus_callarias.htm"

$ws.Range("G5").Value = "3 cubs"
$ws.Range("H5").Value = "https://lionalert.org/page/reproduction-and-offspring"

$ws.Range("G7").Value = " 1 pup"
$ws.Range("H7").Value = "http://www.defenders.org/bats/bats"

$ws.Range("G8").Value = 5
$ws.Range("H8").Value = "https://en.wikipedia.org/wiki/Canine_reproduction"

$ws.Range("G9").Value = 100
$ws.Range("H9").Value = "http://animals.mom.me/reproduction-cycle-dragonfly-9190.html"

$ws.Range("G10").Value = 300
$ws.Range("H10").Value = "http://www.terro.com/pantry-moths-reproduction"

$ws.Range("G11").Value = 10000
$ws.Range("H11").Value = "http://www.ehow.com/about_6465166_do-oak-trees-reproduce_.html"

$ws.Range("G12").Value = 10000

# Update the selected cell to match the saved view state
$ws.Range("C18").Select()
